$wb = $excel.ActiveWorkbook

# Insert the new sheet right after Sheet1, so it becomes Sheet2 / second tab.
$sheet1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Add($null, $sheet1)
$ws2.Name = "Sheet2"

# Populate the new sheet with its content.
$ws2.Range("A1").Value = "additional test sheet"

# Match the recorded selection / active-cell state on Sheet2.
$ws2.Range("F17").Select()

$wb.Save()
